# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6606524410359556, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 0, 4.999471853163143)
    3 = @(3.286832544864788, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 4.840633575959121)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    5 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    # Columns B (2) through G (7)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
